$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Resize the table's columns (widths given in dxa/20 = points)
$t.Columns.Item(1).Width = 24.8
$t.Columns.Item(2).Width = 99.2
$t.Columns.Item(3).Width = 141.75
$t.Columns.Item(4).Width = 134.65
$t.Columns.Item(5).Width = 134.7
$t.Columns.Item(6).Width = 173.4

# Center the text in the first two header cells ("Nr" and "Numer RMA")
$t.Cell(1, 1).Range.ParagraphFormat.Alignment = 1
$t.Cell(1, 2).Range.ParagraphFormat.Alignment = 1
